$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Get-RowData($r) {
    return @{
        A  = $ws.Range("A$r").Value2
        B  = $ws.Range("B$r").Value2
        D  = $ws.Range("D$r").Value2
        E  = $ws.Range("E$r").Value2
        F  = $ws.Range("F$r").Value2
        G  = $ws.Range("G$r").Value2
        H  = $ws.Range("H$r").Value2
        Q  = $ws.Range("Q$r").Value2
        R  = $ws.Range("R$r").Value2
        AC = $ws.Range("AC$r").Value2
    }
}

function Set-RowData($r, $data) {
    $ws.Range("A$r").Value = $data.A
    $ws.Range("B$r").Value = $data.B
    $ws.Range("D$r").Value = $data.D
    $ws.Range("E$r").Value = $data.E
    $ws.Range("F$r").Value = $data.F
    $ws.Range("G$r").Value = $data.G
    $ws.Range("H$r").Value = $data.H
    $ws.Range("Q$r").Value = $data.Q
    $ws.Range("R$r").Value = $data.R
    $ws.Range("AC$r").Value = $data.AC
}

# Snapshot the rows that are involved in the rearrangement before any writes happen.
$row7  = Get-RowData 7
$row8  = Get-RowData 8
$row13 = Get-RowData 13
$row14 = Get-RowData 14
$row15 = Get-RowData 15
$row18 = Get-RowData 18
$row19 = Get-RowData 19

# Rows 7 and 8 swap their species/observation data entirely.
Set-RowData 7 $row8
Set-RowData 8 $row7

# Rows 13, 14, 15 rotate: new13 = old15, new14 = old13, new15 = old14.
Set-RowData 13 $row15
Set-RowData 14 $row13
Set-RowData 15 $row14

# Rows 18 and 19 swap their species/observation data entirely.
Set-RowData 18 $row19
Set-RowData 19 $row18
